$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number + date range) ---
$ws.Range("A8").Value = "Volume 31   Number  50"
$ws.Range("C9").Value = "Report Covering the Week  12/9/2024  Through  12/15/2024"

# --- Cells changing from numeric to the text placeholder "0" (style/shared-string match) ---
# Use Copy from a donor cell that already holds the text "0" in style 13, then restore via the
# same text so the shared string + style are reused exactly like Excel would do.
$ws.Range("D15").Copy($ws.Range("C15"))
$ws.Range("D27").Copy($ws.Range("C27"))
$ws.Range("D22").Copy($ws.Range("C28"))
$ws.Range("D22").Copy($ws.Range("D28"))
$ws.Range("D22").Copy($ws.Range("D29"))
$ws.Range("D22").Copy($ws.Range("D30"))

# --- Cells changing from numeric to the text placeholder "***.*" ---
$ws.Range("E22").Copy($ws.Range("E28"))
$ws.Range("E22").Copy($ws.Range("E29"))
$ws.Range("E22").Copy($ws.Range("E30"))

# --- Row 31: cells changing from text placeholders back to numeric values ---
# First copy format/style from numeric donor cells, then assign the new numeric value.
$ws.Range("I31").Copy($ws.Range("C31"))
$ws.Range("C31").Value = 1
$ws.Range("I31").Copy($ws.Range("D31"))
$ws.Range("D31").Value = 1
$ws.Range("I31").Copy($ws.Range("F31"))
$ws.Range("F31").Value = 1
$ws.Range("I31").Copy($ws.Range("G31"))
$ws.Range("G31").Value = 1
$ws.Range("K31").Copy($ws.Range("E31"))
$ws.Range("E31").Value = 0
$ws.Range("K31").Copy($ws.Range("H31"))
$ws.Range("H31").Value = 0

# --- Plain numeric value updates (no type change) ---
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 200
$ws.Range("C16").Value = 1
$ws.Range("E16").Value = -66.666666666666
$ws.Range("I16").Value = 192
$ws.Range("J16").Value = 199
$ws.Range("K16").Value = -3.517587939698
$ws.Range("L16").Value = -11.111111111111
$ws.Range("M16").Value = -23.809523809523
$ws.Range("N16").Value = -77.304964539007
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 50
$ws.Range("F17").Value = 29
$ws.Range("H17").Value = -9.375
$ws.Range("I17").Value = 380
$ws.Range("J17").Value = 352
$ws.Range("K17").Value = 7.954545454545
$ws.Range("L17").Value = 15.501519756838
$ws.Range("M17").Value = 162.068965517241
$ws.Range("N17").Value = -3.79746835443
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -75
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = -25
$ws.Range("I18").Value = 118
$ws.Range("J18").Value = 136
$ws.Range("K18").Value = -13.235294117647
$ws.Range("L18").Value = -20.27027027027
$ws.Range("M18").Value = -59.726962457337
$ws.Range("N18").Value = -89.948892674616
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = -68.75
$ws.Range("F19").Value = 27
$ws.Range("G19").Value = 43
$ws.Range("H19").Value = -37.209302325581
$ws.Range("I19").Value = 503
$ws.Range("J19").Value = 597
$ws.Range("K19").Value = -15.74539363484
$ws.Range("L19").Value = -22.970903522205
$ws.Range("M19").Value = 43.714285714285
$ws.Range("N19").Value = -11.908931698774
$ws.Range("C20").Value = 8
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 300
$ws.Range("F20").Value = 22
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = 100
$ws.Range("I20").Value = 320
$ws.Range("J20").Value = 244
$ws.Range("K20").Value = 31.147540983606
$ws.Range("L20").Value = 52.380952380952
$ws.Range("M20").Value = 3.559870550161
$ws.Range("N20").Value = -89.918084436042
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 29
$ws.Range("E21").Value = -27.586206896551
$ws.Range("F21").Value = 98
$ws.Range("G21").Value = 108
$ws.Range("H21").Value = -9.259259259259
$ws.Range("I21").Value = 1534
$ws.Range("J21").Value = 1551
$ws.Range("K21").Value = -1.096067053513
$ws.Range("L21").Value = -2.911392405063
$ws.Range("M21").Value = 11.889132020423
$ws.Range("N21").Value = -75.305859626529
$ws.Range("L22").Value = 0
$ws.Range("C24").Value = 32
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = 77.777777777777
$ws.Range("F24").Value = 124
$ws.Range("G24").Value = 81
$ws.Range("H24").Value = 53.086419753086
$ws.Range("I24").Value = 1292
$ws.Range("J24").Value = 1266
$ws.Range("K24").Value = 2.053712480252
$ws.Range("L24").Value = -1.674277016742
$ws.Range("M24").Value = 93.703148425787
$ws.Range("C25").Value = 21
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 200
$ws.Range("F25").Value = 77
$ws.Range("G25").Value = 35
$ws.Range("H25").Value = 120
$ws.Range("I25").Value = 725
$ws.Range("J25").Value = 585
$ws.Range("K25").Value = 23.931623931623
$ws.Range("L25").Value = 27.19298245614
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = -12.5
$ws.Range("F26").Value = 35
$ws.Range("G26").Value = 38
$ws.Range("H26").Value = -7.894736842105
$ws.Range("I26").Value = 642
$ws.Range("J26").Value = 512
$ws.Range("K26").Value = 25.390625
$ws.Range("L26").Value = 26.129666011787
$ws.Range("M26").Value = 35.157894736842
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 3
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 100
$ws.Range("L28").Value = -21.428571428571
$ws.Range("I31").Value = 6
$ws.Range("J31").Value = 5
$ws.Range("K31").Value = 20
$ws.Range("L31").Value = 50
